# Revert "merged master to login branch"
# - Products sheet (sheet1): column E goes back to "supp_ID" (supplier code
#   text like SU001..SU004) instead of the "num_in_stock" numeric column that
#   had been merged in.
# - Sales sheet (sheet4): drop the "num_in_stock" column D that had been
#   merged in, restoring the sheet to 3 data columns (sale_ID, sale_date,
#   prod_ID).
# - Restore view state (active sheet / selections) to what they were before
#   the merge.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Products sheet: restore the "supp_ID" column (was "num_in_stock")
# ---------------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item(1)

# Column E used to carry the bold "row label" style (same as column A's data
# cells, xf index 1) rather than the style used for the num_in_stock values
# (xf index 4). Copy A2's format onto E2:E9 before overwriting the values so
# the cells pick up the correct style.
$wsProducts.Range("A2").Copy() | Out-Null
$wsProducts.Range("E2:E9").PasteSpecial(-4122) | Out-Null

$wsProducts.Range("E1").Value = "supp_ID"
$wsProducts.Range("E2").Value = "SU001"
$wsProducts.Range("E3").Value = "SU001"
$wsProducts.Range("E4").Value = "SU002"
$wsProducts.Range("E5").Value = "SU003"
$wsProducts.Range("E6").Value = "SU003"
$wsProducts.Range("E7").Value = "SU003"
$wsProducts.Range("E8").Value = "SU004"
$wsProducts.Range("E9").Value = "SU004"

$wsProducts.Range("A2:A9").Select() | Out-Null

# ---------------------------------------------------------------------
# Sales sheet: delete the merged-in "num_in_stock" column (D)
# ---------------------------------------------------------------------
$wsSales = $wb.Worksheets.Item(4)
$wsSales.Columns.Item(4).Delete() | Out-Null

$wsSales.Activate() | Out-Null
$wsSales.Range("E25").Select() | Out-Null
